# Team Meet Up 02 Jan 2024
# Adds a new "02-Jan-2024" attendance column (F) to the batch sheet,
# copying Present/Reason status plus the two "Reason" comments, extends
# the data validation range to cover the new columns, and moves the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date header in F1, matching the date format already used by C1:E1
$ws.Range("F1").Value = 45293
$ws.Range("F1").NumberFormat = $ws.Range("E1").NumberFormat

# Attendance status for the new date column
$ws.Range("F2").Value = "Present"
$ws.Range("F3").Value = "Present"
$ws.Range("F4").Value = "Present"
$ws.Range("F5").Value = "Reason"
$ws.Range("F6").Value = "Present"
$ws.Range("F7").Value = "Reason"
$ws.Range("F8").Value = "Present"

# Explanatory comments for the two "Reason" entries
$excel.UserName = "Hp"
$ws.Range("F5").AddComment("Hp:`nHealth Issue")
$ws.Range("F7").AddComment("Hp:`nOutside the Pune")

# Extend the list data validation so it also covers the new columns
$ws.Range("C2:E8").Validation.Delete()
$ws.Range("C2:N8").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# Move the active selection
$ws.Range("G7").Select() | Out-Null
